$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 353 (pushes old rows 353:409 down to 354:410,
# dimension grows from A1:R409 to A1:R410).
$ws.Rows.Item(353).Insert()

# Populate the newly inserted row with the new reading.
$ws.Cells.Item(353, 1).Value = 7
$ws.Cells.Item(353, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(353, 3).Value = "Ñuble"
$ws.Cells.Item(353, 4).Value = 44984
$ws.Cells.Item(353, 5).Value = 16
$ws.Cells.Item(353, 6).Value = 100112002
$ws.Cells.Item(353, 7).Value = "Pimiento"
$ws.Cells.Item(353, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(353, 9).Value = "Primera"
$ws.Cells.Item(353, 10).Value = 60
$ws.Cells.Item(353, 11).Value = 9500
$ws.Cells.Item(353, 12).Value = 10000
$ws.Cells.Item(353, 13).Value = 9750
$ws.Cells.Item(353, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(353, 15).Value = "Región del Maule"
$ws.Cells.Item(353, 16).Value = 542
$ws.Cells.Item(353, 17).Value = 18
$ws.Cells.Item(353, 18).Value = "Hortaliza"
